$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change signer block (rows 41-44, column A) and representative name (E42)
$ws.Range("A41").Value = "Главный бухгалтер"
$ws.Range("A42").Value = "Кахно А.В.,"
$ws.Range("A43").Value = "действующая на основании"
$ws.Range("A44").Value = "Доверенности N 40 от 08.09.2020"

$ws.Range("E42").Value = '${representative_full_name},'
$ws.Range("E43").Value = "действующий на основании"

# Update selection / scroll position to match new target cell
$excel.ActiveWindow.ScrollRow = 30
$ws.Range("C40").Select()
